$p = $ppt.ActivePresentation
$css = $p.ColorSchemes
Write-Output "Count=$($css.Count)"
try {
  $new = $css.Add()
  Write-Output "Added: $new"
} catch {
  Write-Output "ERR add: $_"
}
